$d = $word.ActiveDocument

$replacements = @(
    @("40÷6=6, 4", "25÷9=2, 7"),
    @("84÷2=42, 0", "81÷8=10, 1"),
    @("54÷9=6, 0", "89÷5=17, 4"),
    @("23÷8=2, 7", "58÷5=11, 3"),
    @("39÷6=6, 3", "69÷6=11, 3"),
    @("67÷9=7, 4", "28÷7=4, 0"),
    @("38÷3=12, 2", "25÷7=3, 4"),
    @("61÷8=7, 5", "93÷2=46, 1"),
    @("44÷9=4, 8", "73÷4=18, 1"),
    @("70÷3=23, 1", "33÷8=4, 1"),
    @("90÷5=18, 0", "91÷8=11, 3"),
    @("59÷9=6, 5", "80÷6=13, 2"),
    @("95÷5=19, 0", "11÷5=2, 1"),
    @("31÷6=5, 1", "73÷6=12, 1"),
    @("27÷5=5, 2", "43÷9=4, 7"),
    @("48÷6=8, 0", "66÷5=13, 1"),
    @("40÷5=8, 0", "37÷4=9, 1"),
    @("10÷7=1, 3", "25÷4=6, 1"),
    @("28÷8=3, 4", "63÷5=12, 3"),
    @("95÷2=47, 1", "70÷2=35, 0"),
    @("19÷9=2, 1", "66÷3=22, 0"),
    @("32÷5=6, 2", "79÷7=11, 2"),
    @("39÷2=19, 1", "38÷7=5, 3"),
    @("36÷8=4, 4", "18÷4=4, 2"),
    @("15÷3=5, 0", "21÷7=3, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
